# Auto-generated Excel COM-interop script to apply diff changes
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Assumptions")
$ws2 = $wb.Worksheets.Item("LiabilityCurves")
$ws3 = $wb.Worksheets.Item("Summary")

# --- Assumptions sheet updates ---
$ws1.Range("B3").Value = 0.25
$ws1.Range("D3").Value = 0.2
$ws1.Range("E3").Value = 3
$ws1.Range("F3").Value = 7
$ws1.Range("G3").Value = 15
$ws1.Range("B4").Value = 0.2
$ws1.Range("G4").Value = 17
$ws1.Range("B5").Value = 0.15
$ws1.Range("E5").Value = 3
$ws1.Range("F5").Value = 7
$ws1.Range("G5").Value = 17
$ws1.Range("B6").Value = 0.35
$ws1.Range("D6").Value = 0.2
$ws1.Range("E6").Value = 3
$ws1.Range("G6").Value = 16
$ws1.Range("G7").Value = 17

# --- LiabilityCurves sheet updates ---
$ws2.Range("B2").Value = 0.1689185198207774
$ws2.Range("C2").Value = 0.1589805568300293
$ws2.Range("D2").Value = 0.1292242474524096
$ws2.Range("E2").Value = 0.1446944197272124
$ws2.Range("F2").Value = 0.1663538972060244
$ws2.Range("G2").Value = 0.1433396118284207
$ws2.Range("B3").Value = 0.1349559073074554
$ws2.Range("C3").Value = 0.1286138573237461
$ws2.Range("D3").Value = 0.110208861221367
$ws2.Range("E3").Value = 0.1200092756082437
$ws2.Range("F3").Value = 0.1316541803165323
$ws2.Range("G3").Value = 0.1202472618056763
$ws2.Range("B4").Value = 0.1089267607465874
$ws2.Range("C4").Value = 0.1051107795844779
$ws2.Range("D4").Value = 0.09443969292366754
$ws2.Range("E4").Value = 0.1003753009834837
$ws2.Range("F4").Value = 0.1054980058819011
$ws2.Range("G4").Value = 0.1013434428142841
$ws2.Range("B5").Value = 0.08879741182910135
$ws2.Range("C5").Value = 0.08674421880123459
$ws2.Range("D5").Value = 0.08130208917836107
$ws2.Range("E5").Value = 0.08460673200273505
$ws2.Range("F5").Value = 0.08560509233035483
$ws2.Range("G5").Value = 0.0858043996641694
$ws2.Range("B6").Value = 0.07308268527578786
$ws2.Range("C6").Value = 0.0722484443646008
$ws2.Range("D6").Value = 0.0703057432888993
$ws2.Range("E6").Value = 0.0718219156953592
$ws2.Range("F6").Value = 0.07032554194664559
$ws2.Range("G6").Value = 0.07297682096128102
$ws2.Range("B7").Value = 0.06069410137305786
$ws2.Range("C7").Value = 0.06069201257941569
$ws2.Range("D7").Value = 0.06105852390324381
$ws2.Range("E7").Value = 0.06136132364161986
$ws2.Range("F7").Value = 0.05846298935923304
$ws2.Range("G7").Value = 0.06234150282250712
$ws2.Range("B8").Value = 0.05083062078527152
$ws2.Range("C8").Value = 0.0513862160531242
$ws2.Range("D8").Value = 0.0532459377499143
$ws2.Range("E8").Value = 0.05272824690687104
$ws2.Range("F8").Value = 0.04914766768916227
$ws2.Range("G8").Value = 0.05348489065617643
$ws2.Range("B9").Value = 0.0428999369664698
$ws2.Range("C9").Value = 0.04381906103019465
$ws2.Range("D9").Value = 0.04661499831020754
$ws2.Range("E9").Value = 0.04554581712660226
$ws2.Range("F9").Value = 0.04174526729669332
$ws2.Range("G9").Value = 0.04607677507062253
$ws2.Range("B10").Value = 0.03646171300472617
$ws2.Range("C10").Value = 0.03760759249273037
$ws2.Range("D10").Value = 0.04096154182770435
$ws2.Range("E10").Value = 0.03952579022045105
$ws2.Range("F10").Value = 0.03579146314242987
$ws2.Range("G10").Value = 0.03985279711466438
$ws2.Range("B11").Value = 0.0311865910005001
$ws2.Range("C11").Value = 0.03246341473727911
$ws2.Range("D11").Value = 0.03612024222392923
$ws2.Range("E11").Value = 0.03444581425917648
$ws2.Range("F11").Value = 0.03094485402337349
$ws2.Range("G11").Value = 0.03460071358686885
$ws2.Range("B12").Value = 0.02682654515256866
$ws2.Range("C12").Value = 0.02816770917744057
$ws2.Range("D12").Value = 0.03195674038767736
$ws2.Range("E12").Value = 0.03013282492458573
$ws2.Range("F12").Value = 0.02695311171788654
$ws2.Range("G12").Value = 0.03014960353898046
$ws2.Range("B13").Value = 0.0231933995512811
$ws2.Range("C13").Value = 0.02455309384764711
$ws2.Range("D13").Value = 0.02836143114394794
$ws2.Range("E13").Value = 0.0264508735726873
$ws2.Range("F13").Value = 0.0236286106511671
$ws2.Range("G13").Value = 0.02636137672615365
$ws2.Range("B14").Value = 0.02014322796433338
$ws2.Range("C14").Value = 0.02149041697548803
$ws2.Range("D14").Value = 0.02524455094033466
$ws2.Range("E14").Value = 0.02329216780318401
$ws2.Range("F14").Value = 0.02083086407737565
$ws2.Range("G14").Value = 0.02312408482740998
$ws2.Range("B15").Value = 0.01756499576346825
$ws2.Range("C15").Value = 0.01887911352500871
$ws2.Range("D15").Value = 0.02253228711115453
$ws2.Range("E15").Value = 0.02057044545696172
$ws2.Range("F15").Value = 0.01845384919940954
$ws2.Range("G15").Value = 0.02034664550314149
$ws2.Range("B16").Value = 0.01537226533269883
$ws2.Range("C16").Value = 0.01664013861853028
$ws2.Range("D16").Value = 0.02016369032116088
$ws2.Range("E16").Value = 0.01821604795776315
$ws2.Range("F16").Value = 0.01641684583369911
$ws2.Range("G16").Value = 0.01795467457314591
$ws2.Range("B17").Value = 0.01349711724328469
$ws2.Range("C17").Value = 0.01471076809906965
$ws2.Range("D17").Value = 0.01808821921399954
$ws2.Range("E17").Value = 0.01617223501883976
$ws2.Range("F17").Value = 0.01465780188356605
$ws2.Range("G17").Value = 0.01588718808711703
$ws2.Range("B18").Value = 0.01188567703420053
$ws2.Range("C18").Value = 0.01304075493373511
$ws2.Range("D18").Value = 0.01626378332779362
$ws2.Range("E18").Value = 0.014392409403876
$ws2.Range("F18").Value = 0.01312851751255854
$ws2.Range("G18").Value = 0.01409398794811012
$ws2.Range("B19").Value = 0.01049480801663064
$ws2.Range("C19").Value = 0.01158947271125685
$ws2.Range("D19").Value = 0.01465517927657969
$ws2.Range("E19").Value = 0.01283801161304731
$ws2.Range("F19").Value = 0.01179113968148529
$ws2.Range("G19").Value = 0.01253358525512201
$ws2.Range("B20").Value = 0.009289653063451727
$ws2.Range("C20").Value = 0.01032377994725274
$ws2.Range("D20").Value = 0.01323283781291006
$ws2.Range("E20").Value = 0.01147691006329664
$ws2.Range("F20").Value = 0.01061560196782632
$ws2.Range("G20").Value = 0.01117154716035096
$ws2.Range("B21").Value = 0.008241796424692853
$ws2.Range("C21").Value = 0.009216412601700755
$ws2.Range("D21").Value = 0.01197181707220569
$ws2.Range("E21").Value = 0.0102821597221116
$ws2.Range("F21").Value = 0.009577747331576807
$ws2.Range("G21").Value = 0.009979177743518341
$ws2.Range("B22").Value = 0.007327879960072668
$ws2.Range("C22").Value = 0.008244765253117955
$ws2.Range("D22").Value = 0.01085099113580986
$ws2.Range("E22").Value = 0.009231036384041907
$ws2.Range("F22").Value = 0.008657945198996238
$ws2.Range("G22").Value = 0.008932462712454462
$ws2.Range("B23").Value = 0.006528553778406027
$ws2.Range("C23").Value = 0.007389959585318961
$ws2.Range("D23").Value = 0.009852393881129731
$ws2.Range("E23").Value = 0.008304278543047003
$ws2.Range("F23").Value = 0.007840067126318101
$ws2.Range("G23").Value = 0.008011222832801012
$ws2.Range("B24").Value = 0.005827674125466828
$ws2.Range("C24").Value = 0.006636126401338443
$ws2.Range("D24").Value = 0.008960686571916253
$ws2.Range("E24").Value = 0.007485486763314183
$ws2.Range("F24").Value = 0.007110723271122974
$ws2.Range("G24").Value = 0.007198432795363377
$ws2.Range("B25").Value = 0.005211685056341802
$ws2.Range("C25").Value = 0.00596984728103732
$ws2.Range("D25").Value = 0.008162724292963665
$ws2.Range("E25").Value = 0.006760643490736385
$ws2.Range("F25").Value = 0.00645868916257923
$ws2.Range("G25").Value = 0.006479671468697538
$ws2.Range("B26").Value = 0.004669137540602557
$ws2.Range("C26").Value = 0.005379716393884558
$ws2.Range("D26").Value = 0.007447201551773477
$ws2.Range("E26").Value = 0.006117725744915994
$ws2.Range("F26").Value = 0.005874471849804623
$ws2.Range("G26").Value = 0.005842676719417359
$ws2.Range("B27").Value = 0.004190312026710763
$ws2.Range("C27").Value = 0.004855993409570702
$ws2.Range("D27").Value = 0.006804361467505765
$ws2.Range("E27").Value = 0.005546390070828585
$ws2.Range("F27").Value = 0.005349978589746575
$ws2.Range("G27").Value = 0.005276983652136308
$ws2.Range("B28").Value = 0.003766919462447935
$ws2.Range("C28").Value = 0.00439032602456545
$ws2.Range("D28").Value = 0.006225756188719751
$ws2.Range("E28").Value = 0.005037714217554393
$ws2.Range("F28").Value = 0.004878261367494737
$ws2.Range("G28").Value = 0.0047736295670995
$ws2.Range("B29").Value = 0.003391862283364528
$ws2.Range("C29").Value = 0.003975526138532037
$ws2.Range("D29").Value = 0.005704048716391541
$ws2.Range("E29").Value = 0.004583983758022363
$ws2.Range("F29").Value = 0.004453317837873713
$ws2.Range("G29").Value = 0.004324912423393893
$ws2.Range("B30").Value = 0.003059041625583998
$ws2.Range("C30").Value = 0.003605387742841492
$ws2.Range("D30").Value = 0.005232848306187411
$ws2.Range("E30").Value = 0.0041785146354649
$ws2.Range("F30").Value = 0.00406993453704964
$ws2.Range("G30").Value = 0.003924192337492167
$ws2.Range("B31").Value = 0.002763200484657421
$ws2.Range("C31").Value = 0.003274537535830547
$ws2.Range("D31").Value = 0.004806573200134923
$ws2.Range("E31").Value = 0.003815504683966462
$ws2.Range("F31").Value = 0.003723562010113036
$ws2.Range("G31").Value = 0.003565727803423673

# --- Summary sheet updates ---
$ws3.Range("B2").Value = 5.354648814067438
$ws3.Range("C2").Value = 25.61477230635892
$ws3.Range("D2").Value = 0.1689185198207774
$ws3.Range("E2").Value = 0.5746812849797094
$ws3.Range("B3").Value = 5.603951939757283
$ws3.Range("C3").Value = 27.65811079765091
$ws3.Range("D3").Value = 0.1589805568300293
$ws3.Range("E3").Value = 0.5516978569040888
$ws3.Range("F3").Value = 0.25
$ws3.Range("H3").Value = 0.2
$ws3.Range("I3").Value = 3
$ws3.Range("J3").Value = 7
$ws3.Range("K3").Value = 15
$ws3.Range("B4").Value = 6.357645888860254
$ws3.Range("C4").Value = 33.02367349563603
$ws3.Range("D4").Value = 0.1292242474524096
$ws3.Range("E4").Value = 0.4854806340647045
$ws3.Range("F4").Value = 0.2
$ws3.Range("K4").Value = 17
$ws3.Range("B5").Value = 5.927794806340268
$ws3.Range("C5").Value = 29.76149868479486
$ws3.Range("D5").Value = 0.1446944197272124
$ws3.Range("E5").Value = 0.5215076440170341
$ws3.Range("F5").Value = 0.15
$ws3.Range("I5").Value = 3
$ws3.Range("J5").Value = 7
$ws3.Range("K5").Value = 17
$ws3.Range("B6").Value = 5.563484131954427
$ws3.Range("C6").Value = 28.65111766185504
$ws3.Range("D6").Value = 0.1663538972060244
$ws3.Range("E6").Value = 0.5594367176814583
$ws3.Range("F6").Value = 0.35
$ws3.Range("H6").Value = 0.2
$ws3.Range("I6").Value = 3
$ws3.Range("K6").Value = 16
$ws3.Range("B7").Value = 5.890971612613508
$ws3.Range("C7").Value = 29.02660093741991
$ws3.Range("D7").Value = 0.1433396118284207
$ws3.Range("E7").Value = 0.5237115370738316
$ws3.Range("K7").Value = 17
